$wb = $excel.ActiveWorkbook

# --- Sheet: all_tools (row 3) ---
$ws1 = $wb.Worksheets.Item("all_tools")
$ws1.Range("C3").Value = 87
$ws1.Range("D3").Value = 385
$ws1.Range("F3").Value = -0.1781295403599469
$ws1.Range("G3").Value = 0.01301652681931669
$ws1.Range("H3").Value = -0.2506050579036035
$ws1.Range("I3").Value = 0.01190964758677489

# --- Sheet: checker_framework (row 3) ---
$ws2 = $wb.Worksheets.Item("checker_framework")
$ws2.Range("F3").Value = -0.2163280190361007
$ws2.Range("G3").Value = 0.00681274956405711
$ws2.Range("H3").Value = -0.2723771706469902
$ws2.Range("I3").Value = 0.00611470670329826

# --- Sheet: typestate_checker (row 3) ---
$ws3 = $wb.Worksheets.Item("typestate_checker")
$ws3.Range("C3").Value = 85
$ws3.Range("D3").Value = 318
$ws3.Range("F3").Value = -0.1064800639782623
$ws3.Range("G3").Value = 0.1390306358667875
$ws3.Range("H3").Value = -0.1473585047345789
$ws3.Range("I3").Value = 0.1434514157012983
